$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '57.291.99'
$ws.Range('E2').Value = '  -1.09%  '
$ws.Range('D3').Value = '3.103.62'
$ws.Range('E3').Value = '  -0.12%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '524.54'
$ws.Range('E5').Value = '  -0.22%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '137.19'
$ws.Range('E6').Value = '  -3.71%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.999'
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('D8').Value = '3.098.39'
$ws.Range('E8').Value = '  -0.26%  '
$ws.Range('E9').Value = '  +2.39%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '7.33'
$ws.Range('E10').Value = '  +1.35%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.107'
$ws.Range('E11').Value = '  -0.86%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.402'
$ws.Range('E12').Value = '  +2.43%  '
$ws.Range('D13').Value = '3.628.26'
$ws.Range('E13').Value = '  -0.26%  '
$ws.Range('E14').Value = '  +1.60%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '25.51'
$ws.Range('E15').Value = '  -0.53%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0000163'
$ws.Range('E16').Value = '  -1.40%  '
$ws.Range('D17').Value = '57.331.74'
$ws.Range('E17').Value = '  -1.15%  '
$ws.Range('D18').Value = '3.091.81'
$ws.Range('E18').Value = '  -0.53%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '5.92'
$ws.Range('E19').Value = '  -2.96%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.52'
$ws.Range('E20').Value = '  -2.31%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '7.93'
$ws.Range('E21').Value = '  -1.12%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '349.01'
$ws.Range('E22').Value = '  +1.84%  '
$ws.Range('E23').Value = '  +0.04%  '
$ws.Range('B24').Value = 'Litecoin'
$ws.Range('C24').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '68.06'
$ws.Range('E24').Value = '  +1.10%  '
$ws.Range('B25').Value = 'Polygon'
$ws.Range('C25').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.501'
$ws.Range('E25').Value = '  -2.45%  '
$ws.Range('B26').Value = 'Kaspa'
$ws.Range('C26').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.167'
$ws.Range('E26').Value = '  -1.37%  '
$ws.Range('B27').Value = 'Binance-PegBSC-USD'
$ws.Range('C27').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.00'
$ws.Range('E27').Value = '  -0.05%  '
$ws.Range('B28').Value = 'PEPE'
$ws.Range('C28').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D28').Value = '0.0₃0894'
$ws.Range('E28').Value = '  -3.06%  '
$ws.Range('B29').Value = 'USDe'
$ws.Range('C29').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.00'
$ws.Range('E29').Value = '  +0.13%  '
$ws.Range('B30').Value = 'InternetComputer(DFINITY)'
$ws.Range('C30').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.38'
$ws.Range('E30').Value = '  +1.28%  '
$ws.Range('B31').Value = 'PancakeSwap'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.88'
$ws.Range('E31').Value = '  +0.16%  '
$ws.Range('B32').Value = 'RenderToken'
$ws.Range('C32').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '5.99'
$ws.Range('E32').Value = '  -7.55%  '
$ws.Range('B33').Value = 'EthereumClassic'
$ws.Range('C33').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '20.87'
$ws.Range('E33').Value = '  -1.00%  '
$ws.Range('B34').Value = 'NEARProtocol'
$ws.Range('C34').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.03'
$ws.Range('E34').Value = '  +7.49%  '
$ws.Range('B35').Value = 'Fetch.AI'
$ws.Range('C35').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.16'
$ws.Range('E35').Value = '  -3.97%  '
$ws.Range('B36').Value = 'Monero'
$ws.Range('C36').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '159.22'
$ws.Range('E36').Value = '  +0.70%  '
$ws.Range('B37').Value = 'Aptos'
$ws.Range('C37').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '6.07'
$ws.Range('E37').Value = '  -1.99%  '
$ws.Range('B38').Value = 'EnergySwap'
$ws.Range('C38').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '26.30'
$ws.Range('E38').Value = '  -1.03%  '
$ws.Range('B39').Value = 'ImmutableX'
$ws.Range('C39').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.24'
$ws.Range('E39').Value = '  -1.34%  '
$ws.Range('B40').Value = 'Hedera'
$ws.Range('C40').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0659'
$ws.Range('E40').Value = '  -1.55%  '
$ws.Range('B41').Value = 'Stacks'
$ws.Range('C41').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.59'
$ws.Range('E41').Value = '  +1.84%  '
$ws.Range('B42').Value = 'Filecoin'
$ws.Range('C42').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '4.09'
$ws.Range('E42').Value = '  +1.05%  '
$ws.Range('B43').Value = 'Mantle'
$ws.Range('C43').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.694'
$ws.Range('E43').Value = '  +1.44%  '
$ws.Range('B44').Value = 'Maker'
$ws.Range('C44').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D44').Value = '2.406.87'
$ws.Range('E44').Value = '  +4.97%  '
$ws.Range('B45').Value = 'OKB'
$ws.Range('C45').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '36.62'
$ws.Range('E45').Value = '  -0.57%  '
$ws.Range('B46').Value = 'FirstDigitalUSD'
$ws.Range('C46').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.999'
$ws.Range('E46').Value = '  -0.08%  '
$ws.Range('B47').Value = 'RenzoRestakedETH'
$ws.Range('C47').Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range('D47').Value = '3.136.05'
$ws.Range('E47').Value = '  -0.35%  '
$ws.Range('E48').Value = '  +0.80%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.972'
$ws.Range('E49').Value = '  -3.36%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '5.98'
$ws.Range('E50').Value = '  -1.69%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.767'
$ws.Range('E51').Value = '  +2.69%  '
